# Adds three blank paragraphs and a new paragraph containing a URL after
# the "Get all the shifts for a user, based on " paragraph (the last
# paragraph in the document), moving the _GoBack bookmark from that
# paragraph to the very end of the newly-added URL text.

$d = $word.ActiveDocument

# The _GoBack bookmark currently sits at the end of the last paragraph;
# remove it for now, we'll re-create it in the right spot once the new
# content has been inserted.
$hadGoBack = $d.Bookmarks.Exists("_GoBack")
if ($hadGoBack) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Start from the end of the document (end of the last paragraph, i.e.
# right after "Get all the shifts for a user, based on ").
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)

# Insert three blank paragraphs.
for ($i = 1; $i -le 3; $i++) {
    $r.Text = "`r"
    $r.Collapse(0)
}

# Insert a fourth paragraph break, then put the URL text into it.
$r.Text = "`r"
$r.Collapse(0)

$urlRange = $d.Paragraphs.Last.Range.Duplicate()
$urlRange.Collapse(1)
$urlRange.InsertAfter("http://anantgarg.com/2009/03/13/write-your-own-php-mvc-framework-part-1/")

# Re-create the _GoBack bookmark collapsed at the very end of the new
# text. Adding a bookmark collapsed exactly at the end of the document /
# end of a paragraph is unreliable, so temporarily type a marker
# character after the target point, anchor the bookmark there, then
# remove the marker again - the bookmark stays put.
$urlRange.Collapse(0)
$urlRange.InsertAfter("X")

$bmRange = $urlRange.Duplicate()
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$urlRange.Delete()
